$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 10.65406927831711
$ws.Range("C2").Value = 15.430741855631982
$ws.Range("D2").Value = 7.8425795026476139
$ws.Range("E2").Value = 13.303031899786658

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 12.959001009682119
$ws.Range("C3").Value = 15.234204646145873
$ws.Range("D3").Value = 13.473689143994417
$ws.Range("E3").Value = 15.152048798862683

# Update the selection to match the new state (B1:E3)
$ws.Range("B1:E3").Select()
